$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking Price cells so Excel does not
# auto-convert strings like "1.001" into numbers, which would lose
# the original text formatting (trailing zeros, thousands-dot style, etc.)
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Updated price (D) and volume/percentage (E) values
$ws.Range("D2").Value = "26.569.19"
$ws.Range("E2").Value = "  +0.21%  "
$ws.Range("D3").Value = "1.733.80"
$ws.Range("E3").Value = "  +0.35%  "
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").Value = "245.39"
$ws.Range("E5").Value = "  +0.17%  "
$ws.Range("D6").Value = "1.002"
$ws.Range("E6").Value = "  +0.17%  "
$ws.Range("D7").Value = "0.4950"
$ws.Range("E7").Value = "  +2.93%  "
$ws.Range("D8").Value = "0.2671"
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "0.06238"
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").Value = "1.743.91"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "0.07043"
$ws.Range("E11").Value = "  -1.59%  "
$ws.Range("D12").Value = "15.73"
$ws.Range("E12").Value = "  +0.36%  "
$ws.Range("D13").Value = "4.587"
$ws.Range("E13").Value = "  +1.47%  "
$ws.Range("D14").Value = "0.6093"
$ws.Range("E14").Value = "  -1.41%  "
$ws.Range("D15").Value = "77.99"
$ws.Range("E15").Value = "  +1.13%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.18%  "
$ws.Range("D17").Value = "26.568.47"
$ws.Range("E17").Value = "  +0.18%  "
$ws.Range("D18").Value = "1.003"
$ws.Range("E18").Value = "  +0.37%  "
$ws.Range("D19").Value = "0.000007196"
$ws.Range("E19").Value = "  +3.93%  "
$ws.Range("D20").Value = "11.52"
$ws.Range("E20").Value = "  -1.13%  "
$ws.Range("D21").Value = "1.968.60"
$ws.Range("E21").Value = "  +1.02%  "
$ws.Range("D22").Value = "4.532"
$ws.Range("E22").Value = "  +0.01%  "
$ws.Range("D23").Value = "8.683"
$ws.Range("E23").Value = "  -3.07%  "
$ws.Range("D24").Value = "5.268"
$ws.Range("E24").Value = "  -0.20%  "
$ws.Range("D25").Value = "139.13"
$ws.Range("E25").Value = "  +2.18%  "
$ws.Range("D26").Value = "15.40"
$ws.Range("E26").Value = "  +0.45%  "
$ws.Range("D27").Value = "1.429"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "107.21"
$ws.Range("E28").Value = "  +0.41%  "
$ws.Range("D29").Value = "1.740"
$ws.Range("E29").Value = "  -3.22%  "
$ws.Range("D30").Value = "4.020"
$ws.Range("E30").Value = "  +0.89%  "
$ws.Range("D31").Value = "0.08011"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "3.719"
$ws.Range("E32").Value = "  +0.32%  "
$ws.Range("D33").Value = "0.04575"
$ws.Range("E33").Value = "  +0.29%  "
$ws.Range("D34").Value = "1.001"
$ws.Range("E34").Value = "  +0.15%  "
$ws.Range("D35").Value = "2.620"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "1.009"
$ws.Range("E36").Value = "  +1.80%  "
$ws.Range("D37").Value = "0.6339"
$ws.Range("E37").Value = "  -0.32%  "
$ws.Range("D38").Value = "0.9051"
$ws.Range("D39").Value = "2.031"
$ws.Range("E39").Value = "  -3.26%  "
$ws.Range("D40").Value = "2.412"
$ws.Range("E40").Value = "  -0.40%  "
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("E42").Value = "  +0.36%  "
$ws.Range("D43").Value = "101.20"
$ws.Range("E43").Value = "  -3.71%  "
$ws.Range("D44").Value = "5.459"
$ws.Range("E44").Value = "  -3.20%  "
$ws.Range("D45").Value = "0.3912"
$ws.Range("E45").Value = "  +0.08%  "
$ws.Range("D46").Value = "6.846"
$ws.Range("E46").Value = "  -1.08%  "
$ws.Range("D47").Value = "0.1174"
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").Value = "0.05382"
$ws.Range("E48").Value = "  +0.98%  "
$ws.Range("D49").Value = "30.67"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "7.751"
$ws.Range("E50").Value = "  -1.54%  "
$ws.Range("D51").Value = "1.249"
$ws.Range("E51").Value = "  -1.38%  "
